# Update the "Förändrad" (Changed) date column (C) from 45203 (2023-10-04)
# to 45204 (2023-10-05) for all data rows (2 through 303).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C303").Value = 45204
